# Add a new interview-question row (row 25) to the "Algorithms" section of
# Sheet1: Topic="Algorithms", Question text, and a hyperlinked Comment URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row content -------------------------------------------------
$ws.Range("A25").Value = "Algorithms"
$ws.Range("C25").Value = "What is general form of recurrance solution for dividing functions?"

# Wrap the long URL text in D25, then turn it into a hyperlink (adding the
# hyperlink after WrapText keeps the cell on the built-in "Hyperlink" style).
$ws.Range("D25").WrapText = $true
$ws.Hyperlinks.Add($ws.Cells.Item(25, 4), "https://www.youtube.com/watch?v=OynWkEj0S-s&list=PLDN4rrl48XKpZkf03iYFl-O29szjTrs_O&index=27") | Out-Null

# Match the taller row used for wrapped text.
$ws.Rows.Item(25).RowHeight = 30

# --- Update selection/view to reflect where the author ended up ------
$ws.Range("D26").Select()
